$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 16.95887025742561
$ws.Range("C2").Value = 10.41088652709471
$ws.Range("D2").Value = 6.003140476266533
$ws.Range("E2").Value = 12.0785583340616
$ws.Range("G2").Value = 3.665344553063869
$ws.Range("L2").Value = 9.928025477822979
$ws.Range("M2").Value = 15.58384651004511
$ws.Range("O2").Value = 27.02666271157908
$ws.Range("B3").Value = 16.49152096218613
$ws.Range("C3").Value = 10.16093321222414
$ws.Range("D3").Value = 5.888025489195202
$ws.Range("E3").Value = 12.11669486953835
$ws.Range("G3").Value = 3.668079164905103
$ws.Range("L3").Value = 9.937270524233574
$ws.Range("M3").Value = 15.48997162322161
$ws.Range("O3").Value = 27.08812233584521
$ws.Range("B4").Value = 16.20064209182174
$ws.Range("C4").Value = 10.00307006593457
$ws.Range("D4").Value = 5.818046401158304
$ws.Range("E4").Value = 12.14133492826608
$ws.Range("G4").Value = 3.669846613796764
$ws.Range("L4").Value = 9.944334821351845
$ws.Range("M4").Value = 15.43458540616399
$ws.Range("O4").Value = 27.13352824239347
$ws.Range("B5").Value = 16.08130175576957
$ws.Range("C5").Value = 9.93769474046303
$ws.Range("D5").Value = 5.789748900640655
$ws.Range("E5").Value = 12.15168466417772
$ws.Range("G5").Value = 3.670589164826073
$ws.Range("L5").Value = 9.947562841438069
$ws.Range("M5").Value = 15.41259930134026
$ws.Range("O5").Value = 27.15395178635718
$ws.Range("B6").Value = 16.06144229498842
$ws.Range("C6").Value = 9.926777920844222
$ws.Range("D6").Value = 5.78506464256875
$ws.Range("E6").Value = 12.1534219021264
$ws.Range("G6").Value = 3.670713814029009
$ws.Range("L6").Value = 9.948119953471213
$ws.Range("M6").Value = 15.40898430940885
$ws.Range("O6").Value = 27.15745882535444
$ws.Range("B7").Value = 16.1990356369682
$ws.Range("C7").Value = 10.0021925406812
$ws.Range("D7").Value = 5.817663825250387
$ws.Range("E7").Value = 12.14147325725381
$ws.Range("G7").Value = 3.669856537701905
$ws.Range("L7").Value = 9.944376941078575
$ws.Range("M7").Value = 15.43428650572529
$ws.Range("O7").Value = 27.13379591784147
$ws.Range("B8").Value = 16.79864868245331
$ws.Range("C8").Value = 10.32564831203493
$ws.Range("D8").Value = 5.96332886732718
$ws.Range("E8").Value = 12.09145430410718
$ws.Range("G8").Value = 3.666269148409108
$ws.Range("L8").Value = 9.930925320536538
$ws.Range("M8").Value = 15.55102157569571
$ws.Range("O8").Value = 27.04625740770032
$ws.Range("B9").Value = 17.93590464902631
$ws.Range("C9").Value = 10.92269915690965
$ws.Range("D9").Value = 6.25267762921486
$ws.Range("E9").Value = 12.00303805640752
$ws.Range("G9").Value = 3.659932133339646
$ws.Range("L9").Value = 9.915544222031242
$ws.Range("M9").Value = 15.79700878709638
$ws.Range("O9").Value = 26.93579065806934
$ws.Range("B10").Value = 18.7387790129547
$ws.Range("C10").Value = 11.33570553752606
$ws.Range("D10").Value = 6.465073442141511
$ws.Range("E10").Value = 11.94391554407413
$ws.Range("G10").Value = 3.655696909779362
$ws.Range("L10").Value = 9.910926958910036
$ws.Range("M10").Value = 15.98701211125055
$ws.Range("O10").Value = 26.8923635453943
$ws.Range("B11").Value = 19.09524154760062
$ws.Range("C11").Value = 11.51747989273044
$ws.Range("D11").Value = 6.561164970944019
$ws.Range("E11").Value = 11.91827412022958
$ws.Range("G11").Value = 3.653860490223816
$ws.Range("L11").Value = 9.910271853292098
$ws.Range("M11").Value = 16.0752110773486
$ws.Range("O11").Value = 26.88087652962404
$ws.Range("B12").Value = 19.22883897806523
$ws.Range("C12").Value = 11.58539489429938
$ws.Range("D12").Value = 6.597437691767782
$ws.Range("E12").Value = 11.90874372216811
$ws.Range("G12").Value = 3.653177978624833
$ws.Range("L12").Value = 9.910230983894824
$ws.Range("M12").Value = 16.10884159880627
$ws.Range("O12").Value = 26.87772051525839
$ws.Range("B13").Value = 19.20012996118416
$ws.Range("C13").Value = 11.57080967358911
$ws.Range("D13").Value = 6.589631430209978
$ws.Range("E13").Value = 11.91078829582619
$ws.Range("G13").Value = 3.653324397071767
$ws.Range("L13").Value = 9.910230581161388
$ws.Range("M13").Value = 16.10158872869519
$ws.Range("O13").Value = 26.87834705266787
$ws.Range("B14").Value = 19.10626117103942
$ws.Range("C14").Value = 11.52308594194068
$ws.Range("D14").Value = 6.564151672071499
$ws.Range("E14").Value = 11.91748645686804
$ws.Range("G14").Value = 3.653804081439002
$ws.Range("L14").Value = 9.910264342113676
$ws.Range("M14").Value = 16.07797337062468
$ws.Range("O14").Value = 26.88059293042113
$ws.Range("B15").Value = 19.04857962885062
$ws.Range("C15").Value = 11.49373293688011
$ws.Range("D15").Value = 6.54852847538023
$ws.Range("E15").Value = 11.92161261947366
$ws.Range("G15").Value = 3.654099579917002
$ws.Range("L15").Value = 9.910311985933815
$ws.Range("M15").Value = 16.06353775835694
$ws.Range("O15").Value = 26.88212420356756
$ws.Range("B16").Value = 18.71529583182663
$ws.Range("C16").Value = 11.3236997145402
$ws.Range("D16").Value = 6.458779616826638
$ws.Range("E16").Value = 11.94561642817943
$ws.Range("G16").Value = 3.655818733213917
$ws.Range("L16").Value = 9.910998802863013
$ws.Range("M16").Value = 15.9812817215844
$ws.Range("O16").Value = 26.89328106556858
$ws.Range("B17").Value = 18.50849746820002
$ws.Range("C17").Value = 11.21779588115938
$ws.Range("D17").Value = 6.403557841816314
$ws.Range("E17").Value = 11.96066249934637
$ws.Range("G17").Value = 3.656896431249576
$ws.Range("L17").Value = 9.911789897390607
$ws.Range("M17").Value = 15.93125667549816
$ws.Range("O17").Value = 26.90224688789159
$ws.Range("B18").Value = 18.38873391392925
$ws.Range("C18").Value = 11.15631045414561
$ws.Range("D18").Value = 6.371748220175418
$ws.Range("E18").Value = 11.96943465513128
$ws.Range("G18").Value = 3.65752478936305
$ws.Range("L18").Value = 9.912380995795724
$ws.Range("M18").Value = 15.90265113393608
$ws.Range("O18").Value = 26.90818179271072
$ws.Range("B19").Value = 18.34804747689627
$ws.Range("C19").Value = 11.13539556074089
$ws.Range("D19").Value = 6.360971090563069
$ws.Range("E19").Value = 11.97242505752926
$ws.Range("G19").Value = 3.657739001761275
$ws.Range("L19").Value = 9.912604521690984
$ws.Range("M19").Value = 15.89299524111276
$ws.Range("O19").Value = 26.91032470835973
$ws.Range("B20").Value = 18.53059718151312
$ws.Range("C20").Value = 11.22912907058452
$ws.Range("D20").Value = 6.409441501261254
$ws.Range("E20").Value = 11.95904860729958
$ws.Range("G20").Value = 3.656780829796586
$ws.Range("L20").Value = 9.911691604397802
$ws.Range("M20").Value = 15.9365647482119
$ws.Range("O20").Value = 26.90121190225234
$ws.Range("B21").Value = 19.13387130552157
$ws.Range("C21").Value = 11.53712881042751
$ws.Range("D21").Value = 6.571639118693387
$ws.Range("E21").Value = 11.91551418087915
$ws.Range("G21").Value = 3.653662836954817
$ws.Range("L21").Value = 9.910248807602224
$ws.Range("M21").Value = 16.08490367125913
$ws.Range("O21").Value = 26.87990082551703
$ws.Range("B22").Value = 19.52000987875022
$ws.Range("C22").Value = 11.73305124164846
$ws.Range("D22").Value = 6.676957327503859
$ws.Range("E22").Value = 11.88810752434653
$ws.Range("G22").Value = 3.651700210641755
$ws.Range("L22").Value = 9.910513260979878
$ws.Range("M22").Value = 16.18319042638226
$ws.Range("O22").Value = 26.87293325622405
$ws.Range("B23").Value = 19.31470356106806
$ws.Range("C23").Value = 11.62898807419074
$ws.Range("D23").Value = 6.620822289615341
$ws.Range("E23").Value = 11.90263957151278
$ws.Range("G23").Value = 3.65274084714053
$ws.Range("L23").Value = 9.91026186658833
$ws.Range("M23").Value = 16.13061791963037
$ws.Range("O23").Value = 26.87601369734636
$ws.Range("B24").Value = 18.52060860879269
$ws.Range("C24").Value = 11.2240072011139
$ws.Range("D24").Value = 6.406781688631995
$ws.Range("E24").Value = 11.95977786790338
$ws.Range("G24").Value = 3.656833065873006
$ws.Range("L24").Value = 9.911735618088592
$ws.Range("M24").Value = 15.93416448549117
$ws.Range("O24").Value = 26.90167738884582
$ws.Range("B25").Value = 17.63336852235166
$ws.Range("C25").Value = 10.76548736644003
$ws.Range("D25").Value = 6.174261975883985
$ws.Range("E25").Value = 12.02592775949749
$ws.Range("G25").Value = 3.661572257112323
$ws.Range("L25").Value = 9.918529769423294
$ws.Range("M25").Value = 15.72874936884526
$ws.Range("O25").Value = 26.95907646603519
